$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the typo: "role_favoris" -> "role_favori" in cell A19
$ws.Range("A19").Value = "role_favori"

# Reflect the saved selection/view state from the diff: scrolled so
# row 10 is the top-visible row, with E16 as the active cell.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E16").Select()
